# B1--and-B2-PowerPoint.pptx edit
#   1. Re-style the table on slide 5 to use the built-in table style
#      {58E01270-4AA8-47C4-940C-32E25577ADF0} instead of the custom
#      "Table_0" style {E3D0E7C1-1DB7-4F27-A88F-E8B700179BA4}.
#   2. Swap the deck's colour theme from the "Integral / Red Violet"
#      palette over to the default "Office" palette (the font scheme
#      and format scheme are already identical between the two themes,
#      only the colour scheme differs).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{58E01270-4AA8-47C4-940C-32E25577ADF0}")

# --- 2. Theme colours -------------------------------------------------
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB  = 0        # dk1      000000
$colors.Item(2).RGB  = 16777215 # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388  # dk2      44546A
$colors.Item(4).RGB  = 15132391 # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939 # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501  # accent2  ED7D31
$colors.Item(7).RGB  = 10855845 # accent3  A5A5A5
$colors.Item(8).RGB  = 49407    # accent4  FFC000
$colors.Item(9).RGB  = 12874308 # accent5  4472C4
$colors.Item(10).RGB = 4697456  # accent6  70AD47
$colors.Item(11).RGB = 12673797 # hlink    0563C1
$colors.Item(12).RGB = 7491477  # folHlink 954F72
